# S12/G01: TradingView alert schema and routing design
#
# Updates sprint-tracker rows 99-107 on Sheet1:
#   - Fills in "deviations" (F) / "remarks" (H) narrative that had been
#     left blank while S12/G01 was still pending.
#   - Flips S12/G01 (rows 101-102) from "pending" to "implemented" and
#     rewrites their "pending work" (I) notes now that the TradingView
#     schema/routing work has landed.
#   - Adds matching F/H narrative for the Zerodha adapter tasks (rows
#     103-104) and the config-mapping task (row 107), and refreshes the
#     "pending work" (I) notes for rows 103-107 to describe the next
#     concrete steps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($addr, $text) {
    # Write the value, then force the cell back to the sheet's true
    # default formatting (no wrap / bottom aligned) instead of whatever
    # the column's default (wrapped) style would otherwise hand a brand
    # new cell.
    $cell = $ws.Range($addr)
    $cell.Value = $text
    $cell.WrapText = $false
    $cell.VerticalAlignment = -4107
}

# --- Row 99: S11/G03 TB001 -------------------------------------------------
Set-PlainText "F99" "User scoping for alerts/orders is implemented at the API layer by filtering Queue, Orders, and Analytics endpoints by the current user while still including legacy/global rows (user_id IS NULL); deriving user_id from TradingView payloads is deferred to S12."
Set-PlainText "H99" "Order history and manual queue listings now honour the logged-in user when available, and analytics summary/trades endpoints restrict results to that user plus global trades."

# --- Row 100: S11/G03 TB002 -------------------------------------------------
Set-PlainText "F100" "User scoping for alerts/orders is implemented at the API layer by filtering Queue, Orders, and Analytics endpoints by the current user while still including legacy/global rows (user_id IS NULL); deriving user_id from TradingView payloads is deferred to S12."
Set-PlainText "H100" "Order history and manual queue listings now honour the logged-in user when available, and analytics summary/trades endpoints restrict results to that user plus global trades."

# --- Row 101: S12/G01 TB001 -------------------------------------------------
Set-PlainText "F101" "Introduced st_user_id field and per-user routing directly in the webhook handler rather than keeping S12/G01 as pure documentation; full alert_type and multi-broker fields are deferred to later S12 groups."
$ws.Range("G101").Value = "implemented"
Set-PlainText "H101" "TradingView payloads now carry st_user_id and are mapped to User.username; alerts without a valid st_user_id are ignored with a structured system event, and valid alerts create Alert/Order rows with user_id populated."
$ws.Range("I101").Value = "Extend the normalized alert schema with alert_type and strategy aliasing, and generalise routing beyond Zerodha in S12/G02–G04."

# --- Row 102: S12/G01 TB002 -------------------------------------------------
Set-PlainText "F102" "Introduced st_user_id field and per-user routing directly in the webhook handler rather than keeping S12/G01 as pure documentation; full alert_type and multi-broker fields are deferred to later S12 groups."
$ws.Range("G102").Value = "implemented"
Set-PlainText "H102" "TradingView payloads now carry st_user_id and are mapped to User.username; alerts without a valid st_user_id are ignored with a structured system event, and valid alerts create Alert/Order rows with user_id populated."
$ws.Range("I102").Value = "Extend the normalized alert schema with alert_type and strategy aliasing, and generalise routing beyond Zerodha in S12/G02–G04."

# --- Row 103: S12/G02 TB001 (Zerodha adapter) -------------------------------
$ws.Range("F103").Value = "Zerodha adapter will treat TradingView symbols as broker symbols for NSE/BSE by default, with an overridable mapping table for edge cases; it will also capture a structured alert reason derived from TV fields like strategy.order.comment."
$ws.Range("H103").Value = "Adapter maps TV placeholders (ticker, action, contracts, price) into the normalized schema including symbol_display, broker_symbol/broker_exchange, product, and a human-readable reason field used later in analytics."
$ws.Range("I103").Value = "Implement the Zerodha adapter in code, wire it into the webhook path, and add a column for alert reason so it can be queried and shown in analytics/queue."

# --- Row 104: S12/G02 TB002 (per-user broker/account mapping) --------------
$ws.Range("F104").Value = "Zerodha adapter will treat TradingView symbols as broker symbols for NSE/BSE by default, with an overridable mapping table for edge cases; it will also capture a structured alert reason derived from TV fields like strategy.order.comment."
$ws.Range("H104").Value = "Adapter maps TV placeholders (ticker, action, contracts, price) into the normalized schema including symbol_display, broker_symbol/broker_exchange, product, and a human-readable reason field used later in analytics."
$ws.Range("I104").Value = "Implement the Zerodha adapter in code, wire it into the webhook path, and add a column for alert reason so it can be queried and shown in analytics/queue."

# --- Row 105: S12/G03 TB001 (webhook v2) ------------------------------------
$ws.Range("I105").Value = "Once the Zerodha adapter is in place, update the webhook to use it for all TV alerts and tighten validation so missing core fields (side, qty, symbol, st_user_id) cause explicit rejections."

# --- Row 106: S12/G03 TB002 (webhook v2 tests) ------------------------------
$ws.Range("I106").Value = "Once the Zerodha adapter is in place, update the webhook to use it for all TV alerts and tighten validation so missing core fields (side, qty, symbol, st_user_id) cause explicit rejections."

# --- Row 107: S12/G04 TB001 (config-based mapping) --------------------------
$ws.Range("F107").Value = "Config mapping will cover both symbol translation and per-field extraction (side/qty/price/product/alert_type/reason) so new brokers or internal alert producers can be added without code changes."
$ws.Range("H107").Value = "Design JSON/YAML config that describes how each platform maps into the normalized alert schema, including symbol rules and required/optional fields."
$ws.Range("I107").Value = "Define the config format and a loader/validator, then hook it into the adapter layer so future brokers/platforms reuse the same mapping mechanism."
